# Commit: "Added print function, few conversion functions and empty function finder"
#
# Populates the "mapping" sheet with the Column -> Location lookup table
# (B/C/E columns map to cells C3/C5/C7/C8 used by the print/template
# helpers) and fills in the "template" sheet's label column used by the
# print-out layout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "mapping" sheet: extend the Column/Location table (A2:B2 -> A5:B5)
# Column-A values are entered first, then column-B, to mirror the order
# the labels were actually typed in (also keeps shared-string ids in the
# expected order).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("mapping")

$ws2.Range("A3").Value = "B"
$ws2.Range("A4").Value = "C"
$ws2.Range("A5").Value = "E"

$ws2.Range("B2").Value = "C3"
$ws2.Range("B3").Value = "C5"
$ws2.Range("B4").Value = "C7"
$ws2.Range("B5").Value = "C8"

# ---------------------------------------------------------------------
# "template" sheet: label column (with thin spacer rows in-between) plus
# the column widths used by the print layout.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("template")

$ws3.Columns.Item(1).ColumnWidth = 6.08
$ws3.Columns.Item(2).ColumnWidth = 9.7
$ws3.Columns.Item(3).ColumnWidth = 20.7

$ws3.Range("B3").Value = "Id"
$ws3.Rows.Item(4).RowHeight = 7.5

$ws3.Range("B5").Value = "Name"
$ws3.Rows.Item(6).RowHeight = 7.35

$ws3.Range("B7").Value = "Address"
$ws3.Rows.Item(8).RowHeight = 7.35

$ws3.Range("B9").Value = "Balance"

# ---------------------------------------------------------------------
# Final on-screen selections, matching where each sheet was left.
# ---------------------------------------------------------------------
$ws2.Range("B6").Select()
$ws3.Range("C5").Select()

$ws1 = $wb.Worksheets.Item("data")
$ws1.Activate()
